$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Copy the date formatting (yyyy-mm-dd style) from the cell above so the
# new date cell reuses the same style record instead of creating a new one.
$ws.Range("A25").Copy()
$ws.Range("A26").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Append a new row of data (row 26) following the existing dataset pattern
$ws.Range("A26").Value = 43916
$ws.Range("B26").Value = 1012
$ws.Range("C26").Value = 112
$ws.Range("D26").Value = 12
$ws.Range("E26").Value = 100
$ws.Range("F26").Value = 0

# Move the active selection to F27, as in the edited workbook
$ws.Range("F27").Select()
